# Update the "Förändrad" (Changed) date column (C) for rows 2-12
# from 45183 (2023-09-14) to 45184 (2023-09-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
